$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the XML diff: Coin/Link/Price/Volume(1h) refresh
# for Fri Nov 10 15:24:15 UTC 2023 GitHub Actions run.
# Each entry: column letter -> new text value for that row.
$updates = @{
    2 = @{ 'D'='37.141.84'; 'E'='  -1.42%  ' }
    3 = @{ 'D'='2.089.34'; 'E'='  +2.37%  ' }
    4 = @{ 'E'='  -0.03%  ' }
    5 = @{ 'D'='249.93'; 'E'='  -2.20%  ' }
    6 = @{ 'D'='0.651'; 'E'='  -7.02%  ' }
    7 = @{ 'D'='1.00'; 'E'='  +0.05%  ' }
    8 = @{ 'D'='51.02'; 'E'='  +7.41%  ' }
    9 = @{ 'D'='60.50'; 'E'='  +3.38%  ' }
    10 = @{ 'E'='  -4.47%  ' }
    11 = @{ 'D'='0.0739'; 'E'='  -4.98%  ' }
    12 = @{ 'E'='  +4.51%  ' }
    13 = @{ 'D'='15.16'; 'E'='  -2.83%  ' }
    14 = @{ 'D'='2.361.33'; 'E'='  +1.62%  ' }
    15 = @{ 'D'='0.830'; 'E'='  -2.35%  ' }
    16 = @{ 'D'='2.098.47'; 'E'='  +2.83%  ' }
    17 = @{ 'D'='5.07'; 'E'='  -3.61%  ' }
    18 = @{ 'D'='37.064.68'; 'E'='  -1.60%  ' }
    19 = @{ 'D'='72.09'; 'E'='  -5.05%  ' }
    20 = @{ 'D'='0.0₃0821'; 'E'='  -5.53%  ' }
    21 = @{ 'D'='13.22'; 'E'='  -4.98%  ' }
    22 = @{ 'D'='238.60'; 'E'='  -7.54%  ' }
    23 = @{ 'E'='  -2.92%  ' }
    24 = @{ 'E'='  +0.20%  ' }
    25 = @{ 'D'='2.44'; 'E'='  -3.74%  ' }
    26 = @{ 'D'='169.91'; 'E'='  -0.48%  ' }
    27 = @{ 'D'='9.09'; 'E'='  +0.93%  ' }
    28 = @{ 'D'='20.68'; 'E'='  +0.78%  ' }
    29 = @{ 'D'='2.00'; 'E'='  -7.54%  ' }
    30 = @{ 'D'='0.122'; 'E'='  -6.59%  ' }
    31 = @{ 'D'='1.07'; 'E'='  +18.32%  ' }
    32 = @{ 'B'='Gas'; 'C'='https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'; 'D'='21.89'; 'E'='  +2.60%  ' }
    33 = @{ 'B'='Filecoin'; 'C'='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; 'D'='4.48'; 'E'='  -3.39%  ' }
    34 = @{ 'B'='Hedera'; 'C'='https://coinranking.com/coin/jad286TjB+hedera-hbar'; 'D'='0.0606'; 'E'='  -2.71%  ' }
    35 = @{ 'D'='0.0907'; 'E'='  -1.28%  ' }
    36 = @{ 'E'='  -0.05%  ' }
    37 = @{ 'E'='  +4.28%  ' }
    38 = @{ 'B'='InternetComputer(DFINITY)'; 'C'='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; 'D'='4.06'; 'E'='  -7.86%  ' }
    39 = @{ 'B'='WEMIXToken'; 'C'='https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; 'D'='1.82'; 'E'='  -3.50%  ' }
    40 = @{ 'E'='  -9.67%  ' }
    41 = @{ 'D'='17.65'; 'E'='  +0.24%  ' }
    42 = @{ 'E'='  -3.28%  ' }
    43 = @{ 'E'='  +0.13%  ' }
    44 = @{ 'D'='97.70'; 'E'='  -6.47%  ' }
    45 = @{ 'E'='  -4.07%  ' }
    46 = @{ 'D'='0.0875'; 'E'='  +2.54%  ' }
    47 = @{ 'D'='2.98'; 'E'='  +3.51%  ' }
    48 = @{ 'D'='1.308.32'; 'E'='  -5.62%  ' }
    49 = @{ 'D'='6.90'; 'E'='  +4.34%  ' }
    50 = @{ 'D'='2.282.14'; 'E'='  +3.43%  ' }
    51 = @{ 'D'='2.28'; 'E'='  -7.09%  ' }
}

foreach ($row in $updates.Keys) {
    $rowUpdates = $updates[$row]
    foreach ($col in $rowUpdates.Keys) {
        $cellRef = "$col$row"
        $newValue = $rowUpdates[$col]
        $cell = $ws.Range($cellRef)
        if ($col -eq "D") {
            # Price column holds text that can look numeric (e.g. "1.00", "0.651").
            # Force text storage, matching the original inlineStr cells, then restore
            # the cell style so no stray number-format/style gets introduced.
            $origStyle = $cell.Style
            $cell.NumberFormat = "@"
            $cell.Value = $newValue
            $cell.Style = $origStyle
        } else {
            $cell.Value = $newValue
        }
    }
}
